$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so that numeric-looking
# strings (e.g. "233.20", "0.616") are preserved exactly as text instead of
# being auto-converted to numbers (which would drop formatting/trailing zeros).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Price (column D) updates
$ws.Range("D2").Value = "35.010.90"
$ws.Range("D3").Value = "1.815.20"
$ws.Range("D5").Value = "233.20"
$ws.Range("D6").Value = "0.616"
$ws.Range("D8").Value = "40.32"
$ws.Range("D9").Value = "0.320"
$ws.Range("D10").Value = "0.0685"
$ws.Range("D12").Value = "2.076.20"
$ws.Range("D13").Value = "1.812.34"
$ws.Range("D14").Value = "11.12"
$ws.Range("D16").Value = "0.660"
$ws.Range("D17").Value = "34.995.01"
$ws.Range("D18").Value = "69.40"
$ws.Range("D19").Value = "0.0₃0789"
$ws.Range("D20").Value = "238.24"
$ws.Range("D24").Value = "2.27"
$ws.Range("D25").Value = "172.83"
$ws.Range("D26").Value = "7.83"
$ws.Range("D31").Value = "3.339.07"
$ws.Range("D32").Value = "0.0555"
$ws.Range("D34").Value = "3.97"
$ws.Range("D41").Value = "1.308.34"
$ws.Range("D42").Value = "0.987"
$ws.Range("D45").Value = "2.45"
$ws.Range("D49").Value = "1.991.45"
$ws.Range("D51").Value = "0.0645"

# Row 36/37 swap: Aave <-> TrustWalletToken (coin moved down one rank, new data)
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.14"
$ws.Range("E36").Value = "  +7.88%  "
$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "93.07"
$ws.Range("E37").Value = "  +4.20%  "

# Row 43/44 swap: InjectiveProtocol <-> RenderToken (coin moved down one rank, new data)
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "14.71"
$ws.Range("E44").Value = "  -4.77%  "

# Restore column D style to Normal so the temporary text format does not
# leave a visible/structural style change on the cells.
$priceRange.Style = "Normal"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -12.11%  "
$ws.Range("E9").Value = "  +7.72%  "
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  +4.54%  "
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +3.91%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  +31.96%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E32").Value = "  +6.64%  "
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +5.13%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("E47").Value = "  +7.19%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("E51").Value = "  +5.46%  "
